$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "star" rows (Io / Asgorath) describing the binary star system.
# Order of assignment matters because it controls the order new strings are
# appended to the shared-strings table: Io, Asgorath, Star, /
$ws.Range("A2").Value = "Io"
$ws.Range("A3").Value = "Asgorath"

$ws.Range("D2").Value = "Star"
$ws.Range("D3").Value = "Star"

$ws.Range("B2").Value = "/"
$ws.Range("B3").Value = "/"

# Update the active selection on the sheet.
$ws.Range("B4").Select() | Out-Null

# Switch the page setup to portrait orientation.
$ws.PageSetup.Orientation = 1
